$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cluster / gene symbol labels referenced by the data rows. The sending-
# cluster set now includes the new "Resolving-Mac" cluster alongside the
# pre-existing ECs / FAPs / Inflammatory-Mac / MuSCs clusters.
$ECs = "ECs"
$FAPs = "FAPs"
$InflammatoryMac = "Inflammatory-Mac"
$MuSCs = "MuSCs"
$ResolvingMac = "Resolving-Mac"
$Ligand = "Efna4"
$Receptor = "Epha5"

# Each row: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 numeric measure columns E..T (column order as in the sheet).
$rows = @(
    @($ECs, $Ligand, $Receptor, $MuSCs, 3, 1, 1.974568, 5.923704, 0.5990695552080697, 0.5990695552080698, 3, 1, 0.02440666666666666, 0.07321999999999999, 0.9884709884709885, 0.9884709884709886, 0.04819262298666666, 0.43373360688, 0.592162875399396, 0.5921628753993963),
    @($ECs, $Ligand, $Receptor, $ResolvingMac, 3, 1, 1.974568, 5.923704, 0.5990695552080697, 0.5990695552080698, 1, 0.3333333333333333, 0.0002846666666666667, 0.000854, 0.01152901152901153, 0.01152901152901153, 0.0005620936906666667, 0.005058843216, 0.006906679808673646, 0.006906679808673647),
    @($FAPs, $Ligand, $Receptor, $MuSCs, 3, 1, 1.020259, 3.060777, 0.309539152527049, 0.309539152527049, 3, 1, 0.02440666666666666, 0.07321999999999999, 0.9884709884709885, 0.9884709884709886, 0.02490112132666666, 0.22411009194, 0.3059704720688842, 0.3059704720688842),
    @($FAPs, $Ligand, $Receptor, $ResolvingMac, 3, 1, 1.020259, 3.060777, 0.309539152527049, 0.309539152527049, 1, 0.3333333333333333, 0.0002846666666666667, 0.000854, 0.01152901152901153, 0.01152901152901153, 0.0002904337286666667, 0.002613903558, 0.003568680458164807, 0.003568680458164807),
    @($InflammatoryMac, $Ligand, $Receptor, $MuSCs, 1, 0.3333333333333333, 0.077601, 0.232803, 0.02354357842003994, 0.02354357842003994, 3, 1, 0.02440666666666666, 0.07321999999999999, 0.9884709884709885, 0.9884709884709886, 0.00189398174, 0.01704583566, 0.02327214423300111, 0.02327214423300112),
    @($InflammatoryMac, $Ligand, $Receptor, $ResolvingMac, 1, 0.3333333333333333, 0.077601, 0.232803, 0.02354357842003994, 0.02354357842003994, 1, 0.3333333333333333, 0.0002846666666666667, 0.000854, 0.01152901152901153, 0.01152901152901153, 0.000022090418, 0.000198813762, 0.0002714341870388275, 0.0002714341870388276),
    @($MuSCs, $Ligand, $Receptor, $MuSCs, 3, 1, 0.1837383333333333, 0.551215, 0.05574487261247628, 0.05574487261247628, 3, 1, 0.02440666666666666, 0.07321999999999999, 0.9884709884709885, 0.9884709884709886, 0.004484440255555555, 0.0403599623, 0.05510218933344376, 0.05510218933344377),
    @($MuSCs, $Ligand, $Receptor, $ResolvingMac, 3, 1, 0.1837383333333333, 0.551215, 0.05574487261247628, 0.05574487261247628, 1, 0.3333333333333333, 0.0002846666666666667, 0.000854, 0.01152901152901153, 0.01152901152901153, 0.00005230417888888889, 0.00047073761, 0.0006426832790325181, 0.0006426832790325181),
    @($ResolvingMac, $Ligand, $Receptor, $MuSCs, 1, 0.3333333333333333, 0.03989166666666667, 0.119675, 0.01210284123236505, 0.01210284123236505, 3, 1, 0.02440666666666666, 0.07321999999999999, 0.9884709884709885, 0.9884709884709886, 0.000973622611111111, 0.008762603499999999, 0.01196330743626331, 0.01196330743626331),
    @($ResolvingMac, $Ligand, $Receptor, $ResolvingMac, 1, 0.3333333333333333, 0.03989166666666667, 0.119675, 0.01210284123236505, 0.01210284123236505, 1, 0.3333333333333333, 0.0002846666666666667, 0.000854, 0.01152901152901153, 0.01152901152901153, 0.00001135582777777778, 0.00010220245, 0.0001395337961017327, 0.0001395337961017327)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
